$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(280).Insert()

$ws.Cells.Item(280, 1).Value = 8
$ws.Cells.Item(280, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(280, 3).Value = "Coquimbo"
$ws.Cells.Item(280, 4).Value = 45212
$ws.Cells.Item(280, 5).Value = 4
$ws.Cells.Item(280, 6).Value = 100112001
$ws.Cells.Item(280, 7).Value = "Berenjena"
$ws.Cells.Item(280, 8).Value = "Sin especificar"
$ws.Cells.Item(280, 9).Value = "Primera"
$ws.Cells.Item(280, 10).Value = 500
$ws.Cells.Item(280, 11).Value = 8500
$ws.Cells.Item(280, 12).Value = 9000
$ws.Cells.Item(280, 13).Value = 8750
$ws.Cells.Item(280, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(280, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(280, 16).Value = 175
$ws.Cells.Item(280, 17).Value = 50
$ws.Cells.Item(280, 18).Value = "Hortaliza"
